$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix E23 style (was using stray font/style 12, should match the rest of column E -> style of D23)
# NOTE: only the style/format changes here, the cell's text value is untouched.
$ws.Range("D23").Copy()
$ws.Range("E23").PasteSpecial(-4122)

# Row 24
$ws.Range("A22:G22").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$ws.Range("A24").Value = 43882
$ws.Range("B24").Value = '2PM-7PM'
$ws.Range("C24").Value = 'Team'
$ws.Range("D24").Value = 'Describe project architecture, social context, and interesting issues/pull requests'
$ws.Range("E24").Value = 'All goals'
$ws.Range("F24").Value = 'Our project uses elements of but does not strictly implement a lot of different architectural styles and patterns, making an accurate architectural diagram difficult. On the other hand, the social context for the project is well-documented because of its forum archives and github metrics.'
$ws.Range("G24").Value = 'Good, finished faster than expected'
$ws.Rows(24).RowHeight = 114.9

# Row 25
$ws.Range("A22:G22").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)
$ws.Range("A25").Value = 43888
$ws.Range("B25").Value = '5PM-8PM'
$ws.Range("C25").Value = 'Class'
$ws.Range("D25").Value = 'Attend lecture'
$ws.Range("E25").Value = 'Learned about design patterns'
$ws.Range("F25").Value = 'There are a ton of design patterns out there, and it would be pretty useful to learn more of them. Good planning saves a lot of time.'
$ws.Range("G25").Value = 'Positive'
$ws.Rows(25).RowHeight = 58.5

# Row 26
$ws.Range("A22:G22").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)
$ws.Range("A26").Value = 43891
$ws.Range("B26").Value = '2:00PM-7:30PM'
$ws.Range("C26").Value = 'Team'
$ws.Range("D26").Value = 'Complete part of homework'
$ws.Range("E26").Value = 'All of homework'
$ws.Range("F26").Value = 'Looking for design patterns was not too hard, but checking to make sure they were actually following the design pattern, and distinguishing similar design patterns, was a little difficult. We expected the issue to take a long time as well, but finished faster than expected, although we did end up switching issues in the middle.'
$ws.Range("G26").Value = 'Good'
$ws.Rows(26).RowHeight = 129.1

# Row 27
$ws.Range("A22:G22").Copy()
$ws.Range("A27:G27").PasteSpecial(-4122)
$ws.Range("A27").Value = 43892
$ws.Range("B27").Value = '1PM-3PM'
$ws.Range("C27").Value = 'Self'
$ws.Range("D27").Value = 'Look more at design patterns'
$ws.Range("E27").Value = 'Read more of refactoringguru, watched some videos'
$ws.Range("F27").Value = 'Tried to look more at design patterns and the kinds of problems they solve. Hoping they will come in handy in the future.'
$ws.Range("G27").Value = 'Good, it’s nice to know I have these resources on hand for when I run into a design issue in the future'
$ws.Rows(27).RowHeight = 58.2

# The trailing filler row metadata moves from row 125 down to the very last
# worksheet row (matches the author's saved selection/scroll state).
$ws.Rows(125).Delete()
$ws.Rows(1048576).RowHeight = 12.8

# Final cursor position left by the author after typing the new entries.
$ws.Range("H27").Select()
